$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9363.64609054819
$ws.Range("C2").Value = 8534.76979603661
$ws.Range("E2").Value = 5050.06369943698
$ws.Range("F2").Value = -9.13193768860052

$ws.Range("B3").Value = 9049.15148992294
$ws.Range("C3").Value = 8558.81955328089
$ws.Range("E3").Value = 4853.43268754739
$ws.Range("F3").Value = 183.677176701178

$ws.Range("B4").Value = 8995.19262625859
$ws.Range("C4").Value = 8437.55210894853
$ws.Range("E4").Value = 5054.69751978511
$ws.Range("F4").Value = 187.010401197235

$ws.Range("B5").Value = 8922.31379343847
$ws.Range("C5").Value = 8271.02133863224
$ws.Range("E5").Value = 4991.33258282178
$ws.Range("F5").Value = 177.431413393918

$ws.Range("B6").Value = 8776.90758286408
$ws.Range("C6").Value = 7365.89363135974
$ws.Range("E6").Value = 4875.36272114807
$ws.Range("F6").Value = 134.885681354492

$ws.Range("B7").Value = 2913.75510449073
$ws.Range("C7").Value = 4790.83937824868
$ws.Range("E7").Value = 4175.00370281279
$ws.Range("F7").Value = -1.58987162243883

$ws.Range("C9").Value = 8321.45576783773
$ws.Range("F9").Value = 214.448680452915

$ws.Range("C10").Value = 8430.78076392359
$ws.Range("F10").Value = 219.003888623159

$ws.Range("C11").Value = 8552.7980306891
$ws.Range("F11").Value = 224.087941405055

$ws.Range("C12").Value = 8615.4033451468
$ws.Range("F12").Value = 226.696496174126

$ws.Range("C13").Value = 8094.95237146652
$ws.Range("F13").Value = 205.011038937448

$ws.Range("C14").Value = 5784.04417736491
$ws.Range("F14").Value = 94.0632365526761

$ws.Range("C15").Value = 5674.41507546685
$ws.Range("F15").Value = 89.1618781645881

$wb.Save()
